$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Cells.Item(33, 8).Value = 43.47826
$ws.Cells.Item(33, 9).Value = 43.47826
$ws.Cells.Item(33, 11).Value = 43.47826
$ws.Cells.Item(33, 13).Value = 185.52174
# Row 64
$ws.Cells.Item(64, 8).Value = 2750
$ws.Cells.Item(64, 9).Value = 2750
$ws.Cells.Item(64, 11).Value = 2750
$ws.Cells.Item(64, 13).Value = -2502
# Row 67
$ws.Cells.Item(67, 8).Value = 2750
$ws.Cells.Item(67, 9).Value = 2750
$ws.Cells.Item(67, 11).Value = 2750
$ws.Cells.Item(67, 13).Value = -1892

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 529
$ws.Cells.Item(2, 9).Value = 503
$ws.Cells.Item(2, 10).Value = 685
$ws.Cells.Item(2, 11).Value = 503
$ws.Cells.Item(2, 12).Value = 685
$ws.Cells.Item(2, 13).Value = -390
$ws.Cells.Item(2, 14).Value = -911
# Row 45
$ws.Cells.Item(45, 8).Value = 1416
$ws.Cells.Item(45, 9).Value = 1416
$ws.Cells.Item(45, 10).Value = 0
$ws.Cells.Item(45, 11).Value = 1416
$ws.Cells.Item(45, 12).Value = 0
$ws.Cells.Item(45, 13).Value = -1039
$ws.Cells.Item(45, 14).ClearContents()
# Row 74
$ws.Cells.Item(74, 8).Value = 930.2632
$ws.Cells.Item(74, 9).Value = 920.2222
$ws.Cells.Item(74, 11).Value = 920.2222
$ws.Cells.Item(74, 13).Value = -46.22220000000004
# Row 77
$ws.Cells.Item(77, 8).Value = 930.2632
$ws.Cells.Item(77, 9).Value = 920.2222
$ws.Cells.Item(77, 11).Value = 4601.111
$ws.Cells.Item(77, 13).Value = -233.1109999999999
# Row 116
$ws.Cells.Item(116, 8).Value = 529
$ws.Cells.Item(116, 9).Value = 503
$ws.Cells.Item(116, 10).Value = 685
$ws.Cells.Item(116, 11).Value = 503
$ws.Cells.Item(116, 12).Value = 685
$ws.Cells.Item(116, 13).Value = 1791
$ws.Cells.Item(116, 14).Value = -5273

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 529
$ws.Cells.Item(3, 9).Value = 503
$ws.Cells.Item(3, 10).Value = 685
$ws.Cells.Item(3, 11).Value = 503
$ws.Cells.Item(3, 12).Value = 685
$ws.Cells.Item(3, 13).Value = -389
$ws.Cells.Item(3, 14).Value = -913
# Row 20
$ws.Cells.Item(20, 8).Value = 929.8571
$ws.Cells.Item(20, 9).Value = 982
$ws.Cells.Item(20, 10).Value = 799.5
$ws.Cells.Item(20, 11).Value = 982
$ws.Cells.Item(20, 12).Value = 799.5
$ws.Cells.Item(20, 13).Value = -735
$ws.Cells.Item(20, 14).Value = -1293.5
# Row 80
$ws.Cells.Item(80, 8).Value = 428
$ws.Cells.Item(80, 9).Value = 56
$ws.Cells.Item(80, 10).Value = 800
$ws.Cells.Item(80, 11).Value = 56
$ws.Cells.Item(80, 12).Value = 800
$ws.Cells.Item(80, 13).Value = 942
$ws.Cells.Item(80, 14).Value = -2796
# Row 83
$ws.Cells.Item(83, 8).Value = 428
$ws.Cells.Item(83, 9).Value = 56
$ws.Cells.Item(83, 10).Value = 800
$ws.Cells.Item(83, 11).Value = 280
$ws.Cells.Item(83, 12).Value = 4000
$ws.Cells.Item(83, 13).Value = 4712
$ws.Cells.Item(83, 14).Value = -13984
# Row 86
$ws.Cells.Item(86, 8).Value = 1755.8334
$ws.Cells.Item(86, 9).Value = 1501.5714
$ws.Cells.Item(86, 10).Value = 2645.75
$ws.Cells.Item(86, 11).Value = 1501.5714
$ws.Cells.Item(86, 12).Value = 2645.75
$ws.Cells.Item(86, 13).Value = -378.5714
$ws.Cells.Item(86, 14).Value = -4891.75
# Row 89
$ws.Cells.Item(89, 8).Value = 1755.8334
$ws.Cells.Item(89, 9).Value = 1501.5714
$ws.Cells.Item(89, 10).Value = 2645.75
$ws.Cells.Item(89, 11).Value = 7507.857
$ws.Cells.Item(89, 12).Value = 13228.75
$ws.Cells.Item(89, 13).Value = -1891.857
$ws.Cells.Item(89, 14).Value = -24460.75
# Row 94
$ws.Cells.Item(94, 8).Value = 2388
$ws.Cells.Item(94, 9).Value = 1908.4445
$ws.Cells.Item(94, 10).Value = 2867.5557
$ws.Cells.Item(94, 11).Value = 1908.4445
$ws.Cells.Item(94, 12).Value = 2867.5557
$ws.Cells.Item(94, 13).Value = -1457.4445
$ws.Cells.Item(94, 14).Value = -3769.5557
# Row 99
$ws.Cells.Item(99, 8).Value = 4002.25
$ws.Cells.Item(99, 9).Value = 3769.6667
$ws.Cells.Item(99, 10).Value = 4700
$ws.Cells.Item(99, 11).Value = 3769.6667
$ws.Cells.Item(99, 12).Value = 4700
$ws.Cells.Item(99, 13).Value = -2271.6667
$ws.Cells.Item(99, 14).Value = -7696

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 19187.596
$ws.Cells.Item(31, 9).Value = 13681.125
$ws.Cells.Item(31, 10).Value = 22576.191
$ws.Cells.Item(31, 11).Value = 13681.125
$ws.Cells.Item(31, 12).Value = 22576.191
$ws.Cells.Item(31, 13).Value = -13386.125
$ws.Cells.Item(31, 14).Value = -23166.191
# Row 34
$ws.Cells.Item(34, 8).Value = 19187.596
$ws.Cells.Item(34, 9).Value = 13681.125
$ws.Cells.Item(34, 10).Value = 22576.191
$ws.Cells.Item(34, 11).Value = 13681.125
$ws.Cells.Item(34, 12).Value = 22576.191
$ws.Cells.Item(34, 13).Value = -13479.125
$ws.Cells.Item(34, 14).Value = -22980.191
# Row 86
$ws.Cells.Item(86, 8).Value = 5149.909
$ws.Cells.Item(86, 9).Value = 4970
$ws.Cells.Item(86, 11).Value = 4970
$ws.Cells.Item(86, 13).Value = -3847
# Row 89
$ws.Cells.Item(89, 8).Value = 5149.909
$ws.Cells.Item(89, 9).Value = 4970
$ws.Cells.Item(89, 11).Value = 24850
$ws.Cells.Item(89, 13).Value = -19234
# Row 99
$ws.Cells.Item(99, 8).Value = 1431135.6
$ws.Cells.Item(99, 9).Value = 1001589.8
$ws.Cells.Item(99, 10).Value = 2505000
$ws.Cells.Item(99, 11).Value = 1001589.8
$ws.Cells.Item(99, 12).Value = 2505000
$ws.Cells.Item(99, 13).Value = -1000091.8
$ws.Cells.Item(99, 14).Value = -2507996
# Row 126
$ws.Cells.Item(126, 8).Value = 1431135.6
$ws.Cells.Item(126, 9).Value = 1001589.8
$ws.Cells.Item(126, 10).Value = 2505000
$ws.Cells.Item(126, 11).Value = 3004769.4
$ws.Cells.Item(126, 12).Value = 7515000
$ws.Cells.Item(126, 13).Value = -3002299.4
$ws.Cells.Item(126, 14).Value = -7519940

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Cells.Item(68, 8).Value = 2497.9092
$ws.Cells.Item(68, 10).Value = 2497.9092
$ws.Cells.Item(68, 12).Value = 7493.7276
$ws.Cells.Item(68, 14).Value = -9115.7276
# Row 71
$ws.Cells.Item(71, 8).Value = 2497.9092
$ws.Cells.Item(71, 10).Value = 2497.9092
$ws.Cells.Item(71, 12).Value = 22481.1828
$ws.Cells.Item(71, 14).Value = -30593.1828
# Row 98
$ws.Cells.Item(98, 8).Value = 3999
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 14).ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Cells.Item(122, 8).Value = 1855.0555
$ws.Cells.Item(122, 9).Value = 1803.25
$ws.Cells.Item(122, 10).Value = 2269.5
$ws.Cells.Item(122, 11).Value = 5409.75
$ws.Cells.Item(122, 12).Value = 6808.5
$ws.Cells.Item(122, 13).Value = -2959.75
$ws.Cells.Item(122, 14).Value = -11708.5
# Row 132
$ws.Cells.Item(132, 8).Value = 462.66666
$ws.Cells.Item(132, 9).Value = 462.66666
$ws.Cells.Item(132, 11).Value = 1387.99998
$ws.Cells.Item(132, 13).Value = 1142.00002

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Cells.Item(136, 8).Value = 3502
$ws.Cells.Item(136, 9).Value = 3502
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 10506
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -7956
$ws.Cells.Item(136, 14).ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Cells.Item(136, 8).Value = 21124.25
$ws.Cells.Item(136, 9).Value = 19832.334
$ws.Cells.Item(136, 11).Value = 59497.00199999999
$ws.Cells.Item(136, 13).Value = -56947.00199999999
